# NewYearTweets.pptx: "extra notes and some development files added"
#
# The "Aim" slide and the "Task management" slide swap places: the
# presentation originally shows "Task management" at slide 2 and "Aim"
# at slide 3; after the edit "Aim" comes first (slide 2) and "Task
# management" follows (slide 3). Everything else about those two
# slides (their shapes, bullet text, autofit settings, etc.) travels
# with them unchanged.

$p = $ppt.ActivePresentation

# Locate the two slides by their current title text, rather than
# hard-coding indices, so the script is robust to the starting order.
$aimIndex = $null
$taskIndex = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $title = $p.Slides.Item($i).Shapes.Item(1).TextFrame.TextRange.Text
    if ($title -eq "Aim") {
        $aimIndex = $i
    } elseif ($title -eq "Task management") {
        $taskIndex = $i
    }
}

# Move the "Aim" slide to just before the "Task management" slide.
$aimSlide = $p.Slides.Item($aimIndex)
if ($aimIndex -gt $taskIndex) {
    # Aim currently follows Task management: dropping it in at
    # Task management's current slot pushes Task management one
    # place later, landing Aim immediately ahead of it.
    $aimSlide.MoveTo($taskIndex)
} else {
    # Aim currently precedes Task management: dropping it in right
    # after Task management's slot puts it immediately ahead once
    # Task management shifts up to fill the gap.
    $aimSlide.MoveTo($taskIndex - 1)
}
